# Generate Report for Handoff
# Updates the localization-status report to reflect a newly generated
# handoff round (new source UUID / content hash) and resets the
# per-language handback columns since no handback exists yet for it.

$wb = $excel.ActiveWorkbook

$oldUuid = "a04b0430-406c-4e47-9e35-46b35a874fe4"
$newUuid = "4820b4f5-0d60-42bf-a379-288e2519be1f"
$newHash = "27dd8e1997d244d7cbfcc142bf10a5b7e33a5efc"

$newFileName = "$newUuid.md"
$newPathName = "e2e\$newUuid.md"

# ------------------------------------------------------------------
# Overview sheet
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathName
$wsOverview.Range("G2").Value = "2016-08-17 08:58:46"

# Refresh the B2 hyperlink's display text, keeping the same target.
$overviewHlAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc43797710dabeb9e82def370c1eb0652317f6d4/e2e/" + $oldUuid + ".md"
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewHlAddr, "", "", $newPathName)

# ------------------------------------------------------------------
# zh-cn sheet
# ------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFileName
$wsZhCn.Range("G2").Value = "$newUuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-17 08:58:41"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

$zhcnHlAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc43797710dabeb9e82def370c1eb0652317f6d4/e2e/" + $oldUuid + ".md"
$wsZhCn.Range("A2:I2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhcnHlAddr, "", "", $newFileName)

$wsZhCn.Columns.Item(9).ColumnWidth = $wsZhCn.Columns.Item(9).ColumnWidth
$wsZhCn.Columns.Item(9).AutoFit()
$wsZhCn.Columns.Item(10).AutoFit()

# ------------------------------------------------------------------
# de-de sheet
# ------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFileName
$wsDeDe.Range("G2").Value = "$newUuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-17 08:58:46"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$dedeHlAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc43797710dabeb9e82def370c1eb0652317f6d4/e2e/" + $oldUuid + ".md"
$wsDeDe.Range("A2:I2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $dedeHlAddr, "", "", $newFileName)

$wsDeDe.Columns.Item(9).AutoFit()
$wsDeDe.Columns.Item(10).AutoFit()
